# Merge the split "<id>" / "p103r_N" / "</id>" runs into a single run
# (keeping the formatting of the first "<id>" run: Courier New / 7f6000 / sz18),
# for every occurrence in the document.
$d = $word.ActiveDocument

$searchStart = 0
while ($true) {
    $docEnd = $d.Content.End

    # Find the opening "<id>" tag starting from $searchStart.
    $openRange = $d.Range($searchStart, $docEnd)
    $foundOpen = $openRange.Find.Execute("<id>", $false, $false, $false, $false, $false, `
                                          $true, 1, $false, "", 0)
    if (-not $foundOpen) { break }

    $openStart = $openRange.Start
    $openEnd = $openRange.End

    # Find the matching closing "</id>" tag right after it.
    $closeRange = $d.Range($openEnd, $d.Content.End)
    $foundClose = $closeRange.Find.Execute("</id>", $false, $false, $false, $false, $false, `
                                            $true, 1, $false, "", 0)
    if (-not $foundClose) { break }

    $closeStart = $closeRange.Start
    $closeEnd = $closeRange.End

    # Text of the middle run(s) sitting between "<id>" and "</id>".
    $middleRange = $d.Range($openEnd, $closeStart)
    $middleText = $middleRange.Text

    # Remove everything after the "<id>" run up to (and including) "</id>".
    $tailRange = $d.Range($openEnd, $closeEnd)
    $tailRange.Delete()

    # Re-append the middle text plus the closing tag onto the "<id>" run so
    # the whole "<id>...</id>" string lives in a single run using that
    # run's formatting.
    $idRunRange = $d.Range($openStart, $openEnd)
    $idRunRange.InsertAfter($middleText + "</id>")

    $searchStart = $openStart + ("<id>" + $middleText + "</id>").Length
}
